$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 46 - this shifts the existing rows 46..66 down to 47..67
$ws.Rows.Item(46).Insert()

# Populate the new row 46 with the new weekly data entry
$ws.Range("A46").Value = 6
$ws.Range("B46").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C46").Value = "Metropolitana"
$ws.Range("D46").Value = 44572
$ws.Range("E46").Value = 13
$ws.Range("F46").Value = "Fruta"
$ws.Range("G46").Value = 100101
$ws.Range("H46").Value = "Berries"
$ws.Range("I46").Value = 100101008
$ws.Range("J46").Value = "Mora"
$ws.Range("K46").Value = "Sin especificar"
$ws.Range("L46").Value = "Primera"
$ws.Range("M46").Value = 250
$ws.Range("N46").Value = 6000
$ws.Range("O46").Value = 6000
$ws.Range("P46").Value = 6000
$ws.Range("Q46").Value = "$/bandeja 2 kilos"
$ws.Range("R46").Value = "Provincia de Linares"
$ws.Range("S46").Value = 3000
$ws.Range("T46").Value = 2
